# Append 14 new literature rows (141-154) to the Master sheet,
# matching the IKD_Literature_Master.xlsx diff (GaN CMOS 2026-02-06 update).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 141 ----
$ws.Cells.Item(141, 1).Value = ''
$ws.Cells.Item(141, 2).Value = 'Future of Power and High‐Speed/RF Electronics with Advanced High‐Performance GaN and Heterogeneous Integration with Silicon CMOS'
$ws.Cells.Item(141, 3).Value = 2026
$ws.Cells.Item(141, 4).Value = 'Wiley'
$ws.Cells.Item(141, 5).Value = 'physica status solidi (a)'
$ws.Cells.Item(141, 6).Value = 'Then, Han Wui; Zubair, Ahmad; Bader, Samuel; Koirala, Pratik; Beumer, Michael; Vora, Heli; Golani, Prafful; Radosavljevic, Marko'
$ws.Cells.Item(141, 7).Value = ''
$ws.Cells.Item(141, 8).Value = '10.1002/pssa.202500873'
$ws.Cells.Item(141, 9).Value = 'https://doi.org/10.1002/pssa.202500873'
$ws.Cells.Item(141, 10).Value = 'Journal'
$ws.Cells.Item(141, 11).Value = 'Co-integration'
$ws.Cells.Item(141, 12).Value = 'Experiment'
$ws.Cells.Item(141, 13).Value = 'Gate Stack'
$ws.Cells.Item(141, 14).Value = ''
$ws.Cells.Item(141, 15).Value = ''
$ws.Cells.Item(141, 16).Value = ''
$ws.Cells.Item(141, 17).Value = 'Future of Power and High‐Speed/RF Electronics with Advanced High‐Performance GaN and Heterogeneous Integration with Silicon CMOS'
$ws.Cells.Item(141, 18).Value = 'High'
$c = $ws.Cells.Item(141, 19)
$c.NumberFormat = "@"
$c.Value = '2026-02-06'
$c.Style = "Normal"
$ws.Cells.Item(141, 20).Value = ''
$ws.Rows.Item(141).AutoFit()

# ---- Row 142 ----
$ws.Cells.Item(142, 1).Value = ''
$ws.Cells.Item(142, 2).Value = 'MATLAB/Simulink Studies on Inverter-Based Systems: Dynamic and Voltage Response Analysis'
$ws.Cells.Item(142, 3).Value = 2026
$ws.Cells.Item(142, 4).Value = 'Springer Science and Business Media LLC'
$ws.Cells.Item(142, 5).Value = ''
$ws.Cells.Item(142, 6).Value = 'ENI, FAVOUR'
$ws.Cells.Item(142, 7).Value = ''
$ws.Cells.Item(142, 8).Value = '10.21203/rs.3.rs-8785842/v1'
$ws.Cells.Item(142, 9).Value = 'https://doi.org/10.21203/rs.3.rs-8785842/v1'
$ws.Cells.Item(142, 10).Value = 'Journal'
$ws.Cells.Item(142, 11).Value = 'Inverter'
$ws.Cells.Item(142, 12).Value = 'Experiment'
$ws.Cells.Item(142, 13).Value = 'Contacts'
$ws.Cells.Item(142, 14).Value = ''
$ws.Cells.Item(142, 15).Value = ''
$ws.Cells.Item(142, 16).Value = ''
$ws.Cells.Item(142, 17).Value = 'MATLAB/Simulink Studies on Inverter-Based Systems: Dynamic and Voltage Response Analysis'
$ws.Cells.Item(142, 18).Value = 'High'
$c = $ws.Cells.Item(142, 19)
$c.NumberFormat = "@"
$c.Value = '2026-02-06'
$c.Style = "Normal"
$ws.Cells.Item(142, 20).Value = ''
$ws.Rows.Item(142).AutoFit()

# ---- Row 143 ----
$ws.Cells.Item(143, 1).Value = ''
$ws.Cells.Item(143, 2).Value = 'Dual-Output, Hybrid-Clamped, Quasi-Five-Level Inverter and Its Modulation Strategy'
$ws.Cells.Item(143, 3).Value = 2026
$ws.Cells.Item(143, 4).Value = 'MDPI AG'
$ws.Cells.Item(143, 5).Value = 'Energies'
$ws.Cells.Item(143, 6).Value = 'Wang, Rutian; Wei, Jiahui; Yu, Yang'
$ws.Cells.Item(143, 7).Value = ''
$ws.Cells.Item(143, 8).Value = '10.3390/en19030856'
$ws.Cells.Item(143, 9).Value = 'https://doi.org/10.3390/en19030856'
$ws.Cells.Item(143, 10).Value = 'Journal'
$ws.Cells.Item(143, 11).Value = 'Inverter'
$ws.Cells.Item(143, 12).Value = 'Experiment'
$ws.Cells.Item(143, 13).Value = 'Contacts'
$ws.Cells.Item(143, 14).Value = ''
$ws.Cells.Item(143, 15).Value = ''
$ws.Cells.Item(143, 16).Value = ''
$ws.Cells.Item(143, 17).Value = 'Dual-Output, Hybrid-Clamped, Quasi-Five-Level Inverter and Its Modulation Strategy'
$ws.Cells.Item(143, 18).Value = 'High'
$c = $ws.Cells.Item(143, 19)
$c.NumberFormat = "@"
$c.Value = '2026-02-06'
$c.Style = "Normal"
$ws.Cells.Item(143, 20).Value = ''
$ws.Rows.Item(143).AutoFit()

# ---- Row 144 ----
$ws.Cells.Item(144, 1).Value = ''
$ws.Cells.Item(144, 2).Value = 'MATLAB/Simulink Studies on Inverter-Based Systems: Dynamic and Voltage Response Analysis'
$ws.Cells.Item(144, 3).Value = 2026
$ws.Cells.Item(144, 4).Value = 'Springer Science and Business Media LLC'
$ws.Cells.Item(144, 5).Value = ''
$ws.Cells.Item(144, 6).Value = 'ENI, FAVOUR'
$ws.Cells.Item(144, 7).Value = ''
$ws.Cells.Item(144, 8).Value = '10.21203/rs.3.rs-8785842/v1'
$ws.Cells.Item(144, 9).Value = 'https://doi.org/10.21203/rs.3.rs-8785842/v1'
$ws.Cells.Item(144, 10).Value = 'Journal'
$ws.Cells.Item(144, 11).Value = 'Inverter'
$ws.Cells.Item(144, 12).Value = 'Experiment'
$ws.Cells.Item(144, 13).Value = 'Contacts'
$ws.Cells.Item(144, 14).Value = ''
$ws.Cells.Item(144, 15).Value = ''
$ws.Cells.Item(144, 16).Value = ''
$ws.Cells.Item(144, 17).Value = 'MATLAB/Simulink Studies on Inverter-Based Systems: Dynamic and Voltage Response Analysis'
$ws.Cells.Item(144, 18).Value = 'High'
$c = $ws.Cells.Item(144, 19)
$c.NumberFormat = "@"
$c.Value = '2026-02-06'
$c.Style = "Normal"
$ws.Cells.Item(144, 20).Value = ''
$ws.Rows.Item(144).AutoFit()

# ---- Row 145 ----
$ws.Cells.Item(145, 1).Value = ''
$ws.Cells.Item(145, 2).Value = 'Dual-Output, Hybrid-Clamped, Quasi-Five-Level Inverter and Its Modulation Strategy'
$ws.Cells.Item(145, 3).Value = 2026
$ws.Cells.Item(145, 4).Value = 'MDPI AG'
$ws.Cells.Item(145, 5).Value = 'Energies'
$ws.Cells.Item(145, 6).Value = 'Wang, Rutian; Wei, Jiahui; Yu, Yang'
$ws.Cells.Item(145, 7).Value = ''
$ws.Cells.Item(145, 8).Value = '10.3390/en19030856'
$ws.Cells.Item(145, 9).Value = 'https://doi.org/10.3390/en19030856'
$ws.Cells.Item(145, 10).Value = 'Journal'
$ws.Cells.Item(145, 11).Value = 'Inverter'
$ws.Cells.Item(145, 12).Value = 'Experiment'
$ws.Cells.Item(145, 13).Value = 'Contacts'
$ws.Cells.Item(145, 14).Value = ''
$ws.Cells.Item(145, 15).Value = ''
$ws.Cells.Item(145, 16).Value = ''
$ws.Cells.Item(145, 17).Value = 'Dual-Output, Hybrid-Clamped, Quasi-Five-Level Inverter and Its Modulation Strategy'
$ws.Cells.Item(145, 18).Value = 'High'
$c = $ws.Cells.Item(145, 19)
$c.NumberFormat = "@"
$c.Value = '2026-02-06'
$c.Style = "Normal"
$ws.Cells.Item(145, 20).Value = ''
$ws.Rows.Item(145).AutoFit()

# ---- Row 146 ----
$ws.Cells.Item(146, 1).Value = ''
$ws.Cells.Item(146, 2).Value = 'Future of Power and High‐Speed/RF Electronics with Advanced High‐Performance GaN and Heterogeneous Integration with Silicon CMOS'
$ws.Cells.Item(146, 3).Value = 2026
$ws.Cells.Item(146, 4).Value = 'Wiley'
$ws.Cells.Item(146, 5).Value = 'physica status solidi (a)'
$ws.Cells.Item(146, 6).Value = 'Then, Han Wui; Zubair, Ahmad; Bader, Samuel; Koirala, Pratik; Beumer, Michael; Vora, Heli; Golani, Prafful; Radosavljevic, Marko'
$ws.Cells.Item(146, 7).Value = ''
$ws.Cells.Item(146, 8).Value = '10.1002/pssa.202500873'
$ws.Cells.Item(146, 9).Value = 'https://doi.org/10.1002/pssa.202500873'
$ws.Cells.Item(146, 10).Value = 'Journal'
$ws.Cells.Item(146, 11).Value = 'Co-integration'
$ws.Cells.Item(146, 12).Value = 'Experiment'
$ws.Cells.Item(146, 13).Value = 'Gate Stack'
$ws.Cells.Item(146, 14).Value = ''
$ws.Cells.Item(146, 15).Value = ''
$ws.Cells.Item(146, 16).Value = ''
$ws.Cells.Item(146, 17).Value = 'Future of Power and High‐Speed/RF Electronics with Advanced High‐Performance GaN and Heterogeneous Integration with Silicon CMOS'
$ws.Cells.Item(146, 18).Value = 'High'
$c = $ws.Cells.Item(146, 19)
$c.NumberFormat = "@"
$c.Value = '2026-02-06'
$c.Style = "Normal"
$ws.Cells.Item(146, 20).Value = ''
$ws.Rows.Item(146).AutoFit()

# ---- Row 147 ----
$ws.Cells.Item(147, 1).Value = ''
$ws.Cells.Item(147, 2).Value = 'Modeling and Validation of Junction Temperature Estimation in High-Power SiC MOSFET Inverters for Electric Vehicle Applications'
$ws.Cells.Item(147, 3).Value = 2026
$ws.Cells.Item(147, 4).Value = 'Associacao Brasileira de Eletronica de Potencia SOBRAEP'
$ws.Cells.Item(147, 5).Value = 'Eletrônica de Potência'
$ws.Cells.Item(147, 6).Value = 'Willers, Leonardo R.; Da Silva e Silva, Paulo Henrique Alves; Rocha, Lucas R.; Vieira, Rodrigo Padilha'
$ws.Cells.Item(147, 7).Value = ''
$ws.Cells.Item(147, 8).Value = '10.18618/rep.e202611'
$ws.Cells.Item(147, 9).Value = 'https://doi.org/10.18618/rep.e202611'
$ws.Cells.Item(147, 10).Value = 'Journal'
$ws.Cells.Item(147, 11).Value = 'n-FET'
$ws.Cells.Item(147, 12).Value = 'Experiment'
$ws.Cells.Item(147, 13).Value = 'Gate Stack'
$ws.Cells.Item(147, 14).Value = ''
$ws.Cells.Item(147, 15).Value = ''
$ws.Cells.Item(147, 16).Value = ''
$ws.Cells.Item(147, 17).Value = 'Modeling and Validation of Junction Temperature Estimation in High-Power SiC MOSFET Inverters for Electric Vehicle Applications'
$ws.Cells.Item(147, 18).Value = 'High'
$c = $ws.Cells.Item(147, 19)
$c.NumberFormat = "@"
$c.Value = '2026-02-06'
$c.Style = "Normal"
$ws.Cells.Item(147, 20).Value = ''
$ws.Rows.Item(147).AutoFit()

# ---- Row 148 ----
$ws.Cells.Item(148, 1).Value = ''
$ws.Cells.Item(148, 2).Value = 'E-mode digitally recessed p-NiO tri-junction HEMT with
                    <i>V</i>
                    BR of 2.5 kV'
$ws.Cells.Item(148, 3).Value = 2026
$ws.Cells.Item(148, 4).Value = 'AIP Publishing'
$ws.Cells.Item(148, 5).Value = 'APL Electronic Devices'
$ws.Cells.Item(148, 6).Value = 'Esteghamat, Amirhossein; Rezaei, Mohammad; Fonollosa, Jon Elipe; Zong, Yuan; Boureau, Victor; Ganeeva, Gulnaz; Matioli, Elison'
$ws.Cells.Item(148, 7).Value = ''
$ws.Cells.Item(148, 8).Value = '10.1063/5.0314238'
$ws.Cells.Item(148, 9).Value = 'https://doi.org/10.1063/5.0314238'
$ws.Cells.Item(148, 10).Value = 'Journal'
$ws.Cells.Item(148, 11).Value = 'n-FET'
$ws.Cells.Item(148, 12).Value = 'Experiment'
$ws.Cells.Item(148, 13).Value = 'Gate Stack'
$ws.Cells.Item(148, 14).Value = ''
$ws.Cells.Item(148, 15).Value = ''
$ws.Cells.Item(148, 16).Value = ''
$ws.Cells.Item(148, 17).Value = 'E-mode digitally recessed p-NiO tri-junction HEMT with
                    <i>V</i>
                    BR of 2.5 kV'
$ws.Cells.Item(148, 18).Value = 'High'
$c = $ws.Cells.Item(148, 19)
$c.NumberFormat = "@"
$c.Value = '2026-02-06'
$c.Style = "Normal"
$ws.Cells.Item(148, 20).Value = ''
$ws.Rows.Item(148).AutoFit()

# ---- Row 149 ----
$ws.Cells.Item(149, 1).Value = ''
$ws.Cells.Item(149, 2).Value = 'MATLAB/Simulink Studies on Inverter-Based Systems: Dynamic and Voltage Response Analysis'
$ws.Cells.Item(149, 3).Value = 2026
$ws.Cells.Item(149, 4).Value = 'Springer Science and Business Media LLC'
$ws.Cells.Item(149, 5).Value = ''
$ws.Cells.Item(149, 6).Value = 'ENI, FAVOUR'
$ws.Cells.Item(149, 7).Value = ''
$ws.Cells.Item(149, 8).Value = '10.21203/rs.3.rs-8785842/v1'
$ws.Cells.Item(149, 9).Value = 'https://doi.org/10.21203/rs.3.rs-8785842/v1'
$ws.Cells.Item(149, 10).Value = 'Journal'
$ws.Cells.Item(149, 11).Value = 'Inverter'
$ws.Cells.Item(149, 12).Value = 'Experiment'
$ws.Cells.Item(149, 13).Value = 'Contacts'
$ws.Cells.Item(149, 14).Value = ''
$ws.Cells.Item(149, 15).Value = ''
$ws.Cells.Item(149, 16).Value = ''
$ws.Cells.Item(149, 17).Value = 'MATLAB/Simulink Studies on Inverter-Based Systems: Dynamic and Voltage Response Analysis'
$ws.Cells.Item(149, 18).Value = 'High'
$c = $ws.Cells.Item(149, 19)
$c.NumberFormat = "@"
$c.Value = '2026-02-06'
$c.Style = "Normal"
$ws.Cells.Item(149, 20).Value = ''
$ws.Rows.Item(149).AutoFit()

# ---- Row 150 ----
$ws.Cells.Item(150, 1).Value = ''
$ws.Cells.Item(150, 2).Value = 'Dual-Output, Hybrid-Clamped, Quasi-Five-Level Inverter and Its Modulation Strategy'
$ws.Cells.Item(150, 3).Value = 2026
$ws.Cells.Item(150, 4).Value = 'MDPI AG'
$ws.Cells.Item(150, 5).Value = 'Energies'
$ws.Cells.Item(150, 6).Value = 'Wang, Rutian; Wei, Jiahui; Yu, Yang'
$ws.Cells.Item(150, 7).Value = ''
$ws.Cells.Item(150, 8).Value = '10.3390/en19030856'
$ws.Cells.Item(150, 9).Value = 'https://doi.org/10.3390/en19030856'
$ws.Cells.Item(150, 10).Value = 'Journal'
$ws.Cells.Item(150, 11).Value = 'Inverter'
$ws.Cells.Item(150, 12).Value = 'Experiment'
$ws.Cells.Item(150, 13).Value = 'Contacts'
$ws.Cells.Item(150, 14).Value = ''
$ws.Cells.Item(150, 15).Value = ''
$ws.Cells.Item(150, 16).Value = ''
$ws.Cells.Item(150, 17).Value = 'Dual-Output, Hybrid-Clamped, Quasi-Five-Level Inverter and Its Modulation Strategy'
$ws.Cells.Item(150, 18).Value = 'High'
$c = $ws.Cells.Item(150, 19)
$c.NumberFormat = "@"
$c.Value = '2026-02-06'
$c.Style = "Normal"
$ws.Cells.Item(150, 20).Value = ''
$ws.Rows.Item(150).AutoFit()

# ---- Row 151 ----
$ws.Cells.Item(151, 1).Value = ''
$ws.Cells.Item(151, 2).Value = 'Future of Power and High‐Speed/RF Electronics with Advanced High‐Performance GaN and Heterogeneous Integration with Silicon CMOS'
$ws.Cells.Item(151, 3).Value = 2026
$ws.Cells.Item(151, 4).Value = 'Wiley'
$ws.Cells.Item(151, 5).Value = 'physica status solidi (a)'
$ws.Cells.Item(151, 6).Value = 'Then, Han Wui; Zubair, Ahmad; Bader, Samuel; Koirala, Pratik; Beumer, Michael; Vora, Heli; Golani, Prafful; Radosavljevic, Marko'
$ws.Cells.Item(151, 7).Value = ''
$ws.Cells.Item(151, 8).Value = '10.1002/pssa.202500873'
$ws.Cells.Item(151, 9).Value = 'https://doi.org/10.1002/pssa.202500873'
$ws.Cells.Item(151, 10).Value = 'Journal'
$ws.Cells.Item(151, 11).Value = 'Co-integration'
$ws.Cells.Item(151, 12).Value = 'Experiment'
$ws.Cells.Item(151, 13).Value = 'Gate Stack'
$ws.Cells.Item(151, 14).Value = ''
$ws.Cells.Item(151, 15).Value = ''
$ws.Cells.Item(151, 16).Value = ''
$ws.Cells.Item(151, 17).Value = 'Future of Power and High‐Speed/RF Electronics with Advanced High‐Performance GaN and Heterogeneous Integration with Silicon CMOS'
$ws.Cells.Item(151, 18).Value = 'High'
$c = $ws.Cells.Item(151, 19)
$c.NumberFormat = "@"
$c.Value = '2026-02-06'
$c.Style = "Normal"
$ws.Cells.Item(151, 20).Value = ''
$ws.Rows.Item(151).AutoFit()

# ---- Row 152 ----
$ws.Cells.Item(152, 1).Value = ''
$ws.Cells.Item(152, 2).Value = 'MATLAB/Simulink Studies on Inverter-Based Systems: Dynamic and Voltage Response Analysis'
$ws.Cells.Item(152, 3).Value = 2026
$ws.Cells.Item(152, 4).Value = 'Springer Science and Business Media LLC'
$ws.Cells.Item(152, 5).Value = ''
$ws.Cells.Item(152, 6).Value = 'ENI, FAVOUR'
$ws.Cells.Item(152, 7).Value = ''
$ws.Cells.Item(152, 8).Value = '10.21203/rs.3.rs-8785842/v1'
$ws.Cells.Item(152, 9).Value = 'https://doi.org/10.21203/rs.3.rs-8785842/v1'
$ws.Cells.Item(152, 10).Value = 'Journal'
$ws.Cells.Item(152, 11).Value = 'Inverter'
$ws.Cells.Item(152, 12).Value = 'Experiment'
$ws.Cells.Item(152, 13).Value = 'Contacts'
$ws.Cells.Item(152, 14).Value = ''
$ws.Cells.Item(152, 15).Value = ''
$ws.Cells.Item(152, 16).Value = ''
$ws.Cells.Item(152, 17).Value = 'MATLAB/Simulink Studies on Inverter-Based Systems: Dynamic and Voltage Response Analysis'
$ws.Cells.Item(152, 18).Value = 'High'
$c = $ws.Cells.Item(152, 19)
$c.NumberFormat = "@"
$c.Value = '2026-02-06'
$c.Style = "Normal"
$ws.Cells.Item(152, 20).Value = ''
$ws.Rows.Item(152).AutoFit()

# ---- Row 153 ----
$ws.Cells.Item(153, 1).Value = ''
$ws.Cells.Item(153, 2).Value = 'Dual-Output, Hybrid-Clamped, Quasi-Five-Level Inverter and Its Modulation Strategy'
$ws.Cells.Item(153, 3).Value = 2026
$ws.Cells.Item(153, 4).Value = 'MDPI AG'
$ws.Cells.Item(153, 5).Value = 'Energies'
$ws.Cells.Item(153, 6).Value = 'Wang, Rutian; Wei, Jiahui; Yu, Yang'
$ws.Cells.Item(153, 7).Value = ''
$ws.Cells.Item(153, 8).Value = '10.3390/en19030856'
$ws.Cells.Item(153, 9).Value = 'https://doi.org/10.3390/en19030856'
$ws.Cells.Item(153, 10).Value = 'Journal'
$ws.Cells.Item(153, 11).Value = 'Inverter'
$ws.Cells.Item(153, 12).Value = 'Experiment'
$ws.Cells.Item(153, 13).Value = 'Contacts'
$ws.Cells.Item(153, 14).Value = ''
$ws.Cells.Item(153, 15).Value = ''
$ws.Cells.Item(153, 16).Value = ''
$ws.Cells.Item(153, 17).Value = 'Dual-Output, Hybrid-Clamped, Quasi-Five-Level Inverter and Its Modulation Strategy'
$ws.Cells.Item(153, 18).Value = 'High'
$c = $ws.Cells.Item(153, 19)
$c.NumberFormat = "@"
$c.Value = '2026-02-06'
$c.Style = "Normal"
$ws.Cells.Item(153, 20).Value = ''
$ws.Rows.Item(153).AutoFit()

# ---- Row 154 ----
$ws.Cells.Item(154, 1).Value = ''
$ws.Cells.Item(154, 2).Value = 'Future of Power and High‐Speed/RF Electronics with Advanced High‐Performance GaN and Heterogeneous Integration with Silicon CMOS'
$ws.Cells.Item(154, 3).Value = 2026
$ws.Cells.Item(154, 4).Value = 'Wiley'
$ws.Cells.Item(154, 5).Value = 'physica status solidi (a)'
$ws.Cells.Item(154, 6).Value = 'Then, Han Wui; Zubair, Ahmad; Bader, Samuel; Koirala, Pratik; Beumer, Michael; Vora, Heli; Golani, Prafful; Radosavljevic, Marko'
$ws.Cells.Item(154, 7).Value = ''
$ws.Cells.Item(154, 8).Value = '10.1002/pssa.202500873'
$ws.Cells.Item(154, 9).Value = 'https://doi.org/10.1002/pssa.202500873'
$ws.Cells.Item(154, 10).Value = 'Journal'
$ws.Cells.Item(154, 11).Value = 'Co-integration'
$ws.Cells.Item(154, 12).Value = 'Experiment'
$ws.Cells.Item(154, 13).Value = 'Gate Stack'
$ws.Cells.Item(154, 14).Value = ''
$ws.Cells.Item(154, 15).Value = ''
$ws.Cells.Item(154, 16).Value = ''
$ws.Cells.Item(154, 17).Value = 'Future of Power and High‐Speed/RF Electronics with Advanced High‐Performance GaN and Heterogeneous Integration with Silicon CMOS'
$ws.Cells.Item(154, 18).Value = 'High'
$c = $ws.Cells.Item(154, 19)
$c.NumberFormat = "@"
$c.Value = '2026-02-06'
$c.Style = "Normal"
$ws.Cells.Item(154, 20).Value = ''
$ws.Rows.Item(154).AutoFit()

